# Auto-generated: update market price / profit columns (H:N) for leves across all job sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9:N9").Value = ,@(50, 0, 50, 0, 50, $null, -388)
$ws.Range("H74:N74").Value = ,@(7277484, 7277484, 0, 7277484, 0, -7276548, $null)
$ws.Range("H77:N77").Value = ,@(7277484, 7277484, 0, 36387420, 0, -36382740, $null)
$ws.Range("H100:N100").Value = ,@(2222.2856, 2533.3333, 1989, 2533.3333, 1989, -1992.3333, -3071)
$ws.Range("H113:N113").Value = ,@(2712.75, 2728.8572, 2600, 2728.8572, 2600, 525.1428000000001, -9108)
$ws.Range("H141:N141").Value = ,@(3962.5, 2679.1667, 5887.5, 8037.500100000001, 17662.5, -2857.500100000001, -28022.5)

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45:N45").Value = ,@(1669.6774, 1597.5454, 1846, 1597.5454, 1846, -1220.5454, -2600)
$ws.Range("H61:N61").Value = ,@(1004.8421, 879.4, 1144.2222, 879.4, 1144.2222, -667.4, -1568.2222)
$ws.Range("H97:N97").Value = ,@(874.5, 847.5, 955.5, 847.5, 955.5, -351.5, -1947.5)
$ws.Range("H122:N122").Value = ,@(908.6667, 708, 1109.3334, 2124, 3328.0002, 326, -8228.0002)
$ws.Range("H132:N132").Value = ,@(1642.1786, 1156.0714, 2128.2856, 3468.2142, 6384.8568, -938.2142000000003, -11444.8568)
$ws.Range("H136:N136").Value = ,@(1004.8421, 879.4, 1144.2222, 2638.2, 3432.6666, -88.19999999999982, -8532.6666)

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94:N94").Value = ,@(1011.2222, 1001.4167, 1030.8334, 1001.4167, 1030.8334, -550.4167, -1932.8334)
$ws.Range("H99:N99").Value = ,@(125002424, 125002424, 0, 125002424, 0, -125000926, $null)
$ws.Range("H107:N107").Value = ,@(15280, 2000, 20971.428, 2000, 20971.428, -80, -24811.428)
$ws.Range("H134:N134").Value = ,@(41072.92, 3146.353, 112712, 9439.059000000001, 338136, -6904.059000000001, -343206)

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58:N58").Value = ,@(7002.0586, 973.8, 9513.833000000001, 973.8, 9513.833000000001, -770.8, -9919.833000000001)
$ws.Range("H132:N132").Value = ,@(2466.5, 1882.8235, 3884, 5648.470499999999, 11652, -3118.470499999999, -16712)
$ws.Range("H134:N134").Value = ,@(2222.9048, 1676.3103, 3442.2307, 5028.9309, 10326.6921, -2493.9309, -15396.6921)
$ws.Range("H136:N136").Value = ,@(7002.0586, 973.8, 9513.833000000001, 2921.4, 28541.499, -371.3999999999996, -33641.499)

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113:N113").Value = ,@(712.1429000000001, 834, 644.44446, 2502, 1933.33338, -332, -6273.33338)

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93:N93").Value = ,@(30000, 0, 30000, 0, 30000, $null, -33744)
$ws.Range("H97:N97").Value = ,@(1920, 1988.3846, 1742.2, 1988.3846, 1742.2, -1492.3846, -2734.2)
$ws.Range("H102:N102").Value = ,@(2499.6, 1499, 3166.6667, 1499, 3166.6667, 123, -6410.6667)
$ws.Range("H122:N122").Value = ,@(693929.0600000001, 878459.75, 1939, 2635379.25, 5817, -2632929.25, -10717)
$ws.Range("H126:N126").Value = ,@(2454.6667, 3090.3333, 1819, 9270.999899999999, 5457, -6800.999899999999, -10397)
$ws.Range("H132:N132").Value = ,@(2987.7307, 2638.0557, 3774.5, 7914.1671, 11323.5, -5384.1671, -16383.5)

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7:N7").Value = ,@(2534.875, 2797.8, 2096.6667, 2797.8, 2096.6667, -2685.8, -2320.6667)
$ws.Range("H40:N40").Value = ,@(919588, 1123607.5, 1500, 1123607.5, 1500, -1123471.5, -1772)
$ws.Range("H61:N61").Value = ,@(22378.8, 27598.5, 1500, 27598.5, 1500, -27396.5, -1904)
$ws.Range("H93:N93").Value = ,@(2878.5557, 2001, 3317.3333, 2001, 3317.3333, -753, -5813.3333)
$ws.Range("H100:N100").Value = ,@(12347701, 22224502, 1700, 22224502, 1700, -22223961, -2782)
$ws.Range("H113:N113").Value = ,@(22378.8, 27598.5, 1500, 27598.5, 1500, -25428.5, -5840)
$ws.Range("H122:N122").Value = ,@(8046.722, 10232.708, 3674.75, 30698.124, 11024.25, -28248.124, -15924.25)
$ws.Range("H126:N126").Value = ,@(2534.875, 2797.8, 2096.6667, 8393.400000000001, 6290.000100000001, -5923.400000000001, -11230.0001)
$ws.Range("H132:N132").Value = ,@(2308.7896, 1660.3636, 3200.375, 4981.0908, 9601.125, -2451.0908, -14661.125)

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122:N122").Value = ,@(2208.5715, 2216.818, 2178.3333, 6650.454000000001, 6534.999899999999, -4200.454000000001, -11434.9999)
$ws.Range("H126:N126").Value = ,@(810.7692, 818.3333, 793.75, 2454.9999, 2381.25, 15.0001000000002, -7321.25)
$ws.Range("H132:N132").Value = ,@(2166.0833, 5504, 1862.6364, 16512, 5587.9092, -13982, -10647.9092)
$ws.Range("H136:N136").Value = ,@(1400.3334, 1329.0714, 1581.7273, 3987.2142, 4745.1819, -1437.2142, -9845.1819)

Write-Host "Applied all Bahamut profit updates"
